# Build site regeneration edit for LOQ4022.xlsx
# Summary of change (per XML diff):
#  - Row 13 (which only held "Docentes responsáveis:" value in B/C, no A label)
#    is removed entirely, shifting all subsequent rows up by one.
#  - A handful of cells end up carrying "stale" / shifted text values as a
#    side effect of the regeneration (this mirrors the shared-strings diff
#    exactly): Objetivos/Programa resumido/Programa/Método/Critério/Norma de
#    recuperação/Bibliografia rows end up showing content that used to
#    belong to a different row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete row 13 entirely -- everything below shifts up one row.
$ws.Rows.Item(13).Delete()

# 2) Fix up the cell values that differ from a plain "shift up" of the
#    old content (these are the actual content changes the diff encodes).

# Row 10 "Objetivos:" -> B/C now show the docente's name instead of the
# original objectives paragraph.
$ws.Range("B10").Value = "4808662 - Lucrécio Fábio dos Santos"
$ws.Range("C10").Value = "4808662 - Lucrécio Fábio dos Santos"

# Row 13 "Programa resumido:" -> B/C now just say "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 "Programa:" -> B/C now show the activation date. "01/01/2022" looks
# like a date to Excel's input parser, so go through a formula + paste-values
# round-trip to force it to land as plain text (matching the original cell
# style/format, no new number format).
$c = $ws.Cells.Item(15, 2)
$c.Formula = "=""01/01/2022"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null

$c = $ws.Cells.Item(15, 3)
$c.Formula = "=""01/01/2022"""
$c.Copy() | Out-Null
$c.PasteSpecial(-4163) | Out-Null

# Row 18 "Método:" -> B/C now show the docente's name.
$ws.Range("B18").Value = "4808662 - Lucrécio Fábio dos Santos"
$ws.Range("C18").Value = "4808662 - Lucrécio Fábio dos Santos"

# Row 19 "Critério:" -> B/C now show the old "Método" description text.
$ws.Range("B19").Value = "Aulas expositivas, desenvolvimento de trabalhos e exercícios em sala e fora de sala de aula, discussão de casos práticos."
$ws.Range("C19").Value = "Aulas expositivas, desenvolvimento de trabalhos e exercícios em sala e fora de sala de aula, discussão de casos práticos."

# Row 20 "Norma de recuperação:" -> B/C now show the old "Critério" text.
$ws.Range("B20").Value = "Provas em sala, entrega de trabalhos e exercícios ou casos práticos elaborados fora de sala de aula."
$ws.Range("C20").Value = "Provas em sala, entrega de trabalhos e exercícios ou casos práticos elaborados fora de sala de aula."

# Row 21 "Bibliografia:" -> B/C now show the old "Norma de recuperação" text.
$ws.Range("B21").Value = "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação."
$ws.Range("C21").Value = "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação."
